$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Apply cell formatting first (copy/paste-special formats only, does not
# touch the shared-string table), mirroring the style families already used
# elsewhere in the sheet:
#   - rows 5 / 7 / 10 / 14 use the "border" family (style 6 for A/B, 7 for
#     C/D/E) and are reused for the new rows 16 and 18.
#   - rows 2 / 6 / 9 / 11 / 12 / 15 use the "header" family (style 4 for
#     A/B, 5 for C/D/E) and are reused for the new rows 17, 19 and 20.
# ---------------------------------------------------------------------------
$ws.Range("A5:E5").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)

$ws.Range("A15:E15").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)

$ws.Range("A5:E5").Copy()
$ws.Range("A18:E18").PasteSpecial(-4122)

$ws.Range("A15:E15").Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)

$ws.Range("A15:E15").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Fill in the cell values, in the exact order the new strings were typed in
# so the shared-string table indices line up with the source workbook.
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = 'SCRIPT/T01P01A/um2504.ssb'

$ws.Range("C17").Value = ' The grand master of all\nthings bad?'
$ws.Range("C18").Value = ' Ummm... Sorry.[K] I\''ve never heard\nof such a creature.'

$ws.Range("A17").Value = 'SCRIPT/T01P01A/us0103.ssb'
$ws.Range("D17").Value = ' Гранд мастер всего самого\nплохого?'
$ws.Range("D18").Value = ' Эммм... Простите.[K] Я ничего не\nзнаю об этом существе.'

$ws.Range("E17").Value = ' Ãñàîä íàòóåñ âòåãï òàíïãï\nðìïöïãï?'
$ws.Range("E18").Value = ' Üííí... Ðñïòóéóå.[K] Ÿ îéœåãï îå\nèîàý ïá üóïí òôþåòóâå.'

$ws.Range("C19").Value = ' We\''ve heard! You graduated\nfrom the guild!'
$ws.Range("C20").Value = ' Congratulations!'

$ws.Range("A19").Value = 'SCRIPT/T01P01A/us0108.ssb'
$ws.Range("D19").Value = ' Мы уже знаем! Вы выпустились\nиз гильдии!'
$ws.Range("D20").Value = ' Поздравляю!'

$ws.Range("E19").Value = ' Íú ôçå èîàåí! Âú âúðôòóéìéòû\néè ãéìûäéé!'
$ws.Range("E20").Value = ' Ðïèäñàâìÿý!'

$ws.Range("A20").Value = 'SCRIPT/T01P01A/us3101.ssb'

# ---------------------------------------------------------------------------
# Numeric "line number" cells.
# ---------------------------------------------------------------------------
$ws.Range("B17").Value = 149
$ws.Range("B18").Value = 152
$ws.Range("B19").Value = 126
$ws.Range("B20").Value = 130

# ---------------------------------------------------------------------------
# Row height tweaks.
# ---------------------------------------------------------------------------
$ws.Rows("15").RowHeight = 48
$ws.Rows("16").RowHeight = 43.2
$ws.Rows("17").RowHeight = 43.2
$ws.Rows("18").RowHeight = 21.6
$ws.Rows("19").RowHeight = 43.2
$ws.Rows("20").RowHeight = 43.2

# ---------------------------------------------------------------------------
# Selection / view state.
# ---------------------------------------------------------------------------
$ws.Range("E20").Select()

Write-Host "done"
